# Auto-generated Excel COM-interop script to apply scheduled price-update edits
# across the eight Leve-profit worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 2650
$ws.Range("I26").Value = 1800
$ws.Range("J26").Value = 3500
$ws.Range("K26").Value = 1800
$ws.Range("L26").Value = 3500
$ws.Range("M26").Value = -1456
$ws.Range("N26").Value = -4188
$ws.Range("H28").Value = 2720.923
$ws.Range("I28").Value = 2657.5
$ws.Range("J28").Value = 2932.3333
$ws.Range("K28").Value = 2657.5
$ws.Range("L28").Value = 2932.3333
$ws.Range("M28").Value = -2172.5
$ws.Range("N28").Value = -3902.3333
$ws.Range("H55").Value = 349.56522
$ws.Range("I55").Value = 464.7
$ws.Range("K55").Value = 464.7
$ws.Range("M55").Value = -250.7
$ws.Range("H74").Value = 6297
$ws.Range("I74").Value = 5373.75
$ws.Range("K74").Value = 5373.75
$ws.Range("M74").Value = -4437.75
$ws.Range("H77").Value = 6297
$ws.Range("I77").Value = 5373.75
$ws.Range("K77").Value = 26868.75
$ws.Range("M77").Value = -22188.75
$ws.Range("H137").Value = 5801.7085
$ws.Range("I137").Value = 2939.7144
$ws.Range("K137").Value = 8819.143199999999
$ws.Range("M137").Value = -6269.143199999999
$ws.Range("H141").Value = 4093.0386
$ws.Range("I141").Value = 2559.4736
$ws.Range("K141").Value = 7678.4208
$ws.Range("M141").Value = -2498.4208

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2812.08
$ws.Range("I61").Value = 2057.6924
$ws.Range("K61").Value = 2057.6924
$ws.Range("M61").Value = -1845.6924
$ws.Range("H74").Value = 2146.7058
$ws.Range("I74").Value = 1588.5
$ws.Range("J74").Value = 3486.4
$ws.Range("K74").Value = 1588.5
$ws.Range("L74").Value = 3486.4
$ws.Range("M74").Value = -714.5
$ws.Range("N74").Value = -5234.4
$ws.Range("H77").Value = 2146.7058
$ws.Range("I77").Value = 1588.5
$ws.Range("J77").Value = 3486.4
$ws.Range("K77").Value = 7942.5
$ws.Range("L77").Value = 17432
$ws.Range("M77").Value = -3574.5
$ws.Range("N77").Value = -26168
$ws.Range("H132").Value = 3509.9092
$ws.Range("J132").Value = 3000
$ws.Range("L132").Value = 9000
$ws.Range("N132").Value = -14060
$ws.Range("H136").Value = 2812.08
$ws.Range("I136").Value = 2057.6924
$ws.Range("K136").Value = 6173.0772
$ws.Range("M136").Value = -3623.0772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1537.762
$ws.Range("I86").Value = 1487.8823
$ws.Range("K86").Value = 1487.8823
$ws.Range("M86").Value = -364.8823
$ws.Range("H89").Value = 1537.762
$ws.Range("I89").Value = 1487.8823
$ws.Range("K89").Value = 7439.4115
$ws.Range("M89").Value = -1823.4115
$ws.Range("H134").Value = 2627.889
$ws.Range("I134").Value = 2579.625
$ws.Range("K134").Value = 7738.875
$ws.Range("M134").Value = -5203.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 188.3077
$ws.Range("I7").Value = 207.5
$ws.Range("K7").Value = 207.5
$ws.Range("M7").Value = -94.5
$ws.Range("H31").Value = 7649.3794
$ws.Range("I31").Value = 5064.75
$ws.Range("K31").Value = 5064.75
$ws.Range("M31").Value = -4769.75
$ws.Range("H34").Value = 7649.3794
$ws.Range("I34").Value = 5064.75
$ws.Range("K34").Value = 5064.75
$ws.Range("M34").Value = -4862.75
$ws.Range("H62").Value = 15795.417
$ws.Range("I62").Value = 14237.692
$ws.Range("J62").Value = 17636.363
$ws.Range("K62").Value = 14237.692
$ws.Range("L62").Value = 17636.363
$ws.Range("M62").Value = -13613.692
$ws.Range("N62").Value = -18884.363
$ws.Range("H65").Value = 15795.417
$ws.Range("I65").Value = 14237.692
$ws.Range("J65").Value = 17636.363
$ws.Range("K65").Value = 71188.45999999999
$ws.Range("L65").Value = 88181.815
$ws.Range("M65").Value = -68068.45999999999
$ws.Range("N65").Value = -94421.815
$ws.Range("H111").Value = 77500
$ws.Range("J111").Value = 77500
$ws.Range("L111").Value = 77500
$ws.Range("N111").Value = -85680

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 9326.333000000001
$ws.Range("I56").Value = 9326.333000000001
$ws.Range("K56").Value = 9326.333000000001
$ws.Range("M56").Value = -8796.333000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9824.643
$ws.Range("I70").Value = 9294.333000000001
$ws.Range("K70").Value = 9294.333000000001
$ws.Range("M70").Value = -9024.333000000001
$ws.Range("H73").Value = 9824.643
$ws.Range("I73").Value = 9294.333000000001
$ws.Range("K73").Value = 9294.333000000001
$ws.Range("M73").Value = -8358.333000000001
$ws.Range("H80").Value = 10428.458
$ws.Range("I80").Value = 3328.3572
$ws.Range("K80").Value = 3328.3572
$ws.Range("M80").Value = -2330.3572
$ws.Range("H83").Value = 10428.458
$ws.Range("I83").Value = 3328.3572
$ws.Range("K83").Value = 16641.786
$ws.Range("M83").Value = -11649.786
$ws.Range("H122").Value = 1915.8334
$ws.Range("I122").Value = 1923.75
$ws.Range("K122").Value = 5771.25
$ws.Range("M122").Value = -3321.25
$ws.Range("H132").Value = 4447.8335
$ws.Range("I132").Value = 3857.4
$ws.Range("K132").Value = 11572.2
$ws.Range("M132").Value = -9042.200000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1736.5
$ws.Range("J22").Value = 2368.6667
$ws.Range("L22").Value = 2368.6667
$ws.Range("N22").Value = -2958.6667
$ws.Range("H27").Value = 1736.5
$ws.Range("J27").Value = 2368.6667
$ws.Range("L27").Value = 2368.6667
$ws.Range("N27").Value = -2582.6667
$ws.Range("H40").Value = 3972.9048
$ws.Range("J40").Value = 3726.5715
$ws.Range("L40").Value = 3726.5715
$ws.Range("N40").Value = -3998.5715
$ws.Range("H93").Value = 2000.65
$ws.Range("I93").Value = 2044.4286
$ws.Range("J93").Value = 1898.5
$ws.Range("K93").Value = 2044.4286
$ws.Range("L93").Value = 1898.5
$ws.Range("M93").Value = -796.4286
$ws.Range("N93").Value = -4394.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H42").Value = 9044
$ws.Range("I42").Value = 9044
$ws.Range("K42").Value = 9044
$ws.Range("M42").Value = -8666
$ws.Range("H113").Value = 476.76923
$ws.Range("I113").Value = 260.1
$ws.Range("J113").Value = 1199
$ws.Range("K113").Value = 780.3000000000001
$ws.Range("L113").Value = 3597
$ws.Range("M113").Value = 1389.7
$ws.Range("N113").Value = -7937
$ws.Range("H132").Value = 2867.2856
$ws.Range("I132").Value = 2843.4473
$ws.Range("K132").Value = 8530.341899999999
$ws.Range("M132").Value = -6000.341899999999
